$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (shifts old D:K data to F:M)
$ws.Columns("D:E").Insert()

# Copy cell formatting (number formats / styles) from the columns that now
# hold the old D:E data (shifted to F:G) onto the freshly inserted D:E so the
# new columns pick up the same date / number styling as their neighbours.
$ws.Range("F7:G102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the two new columns with the latest two quarters of data.
$ws.Range("D7").Value2 = 43464
$ws.Range("E7").Value2 = 43352
$ws.Range("D8").Value2 = 1082100
$ws.Range("E8").Value2 = 786000
$ws.Range("D9").Value2 = 668200
$ws.Range("E9").Value2 = 490700
$ws.Range("D10").Value2 = 413900
$ws.Range("E10").Value2 = 295300
$ws.Range("D12").Value2 = "NA"
$ws.Range("E12").Value2 = "NA"
$ws.Range("D13").Value2 = 0
$ws.Range("E13").Value2 = 0
$ws.Range("D14").Value2 = 0
$ws.Range("E14").Value2 = 0
$ws.Range("D15").Value2 = 0
$ws.Range("E15").Value2 = 0
$ws.Range("D17").Value2 = 902500
$ws.Range("E17").Value2 = 653500
$ws.Range("D18").Value2 = 179600
$ws.Range("E18").Value2 = 132500
$ws.Range("D20").Value2 = 900
$ws.Range("E20").Value2 = 800
$ws.Range("D21").Value2 = 198400
$ws.Range("E21").Value2 = 145700
$ws.Range("D22").Value2 = 46000
$ws.Range("E22").Value2 = 34000
$ws.Range("D23").Value2 = 134600
$ws.Range("E23").Value2 = 99200
$ws.Range("D24").Value2 = 22900
$ws.Range("E24").Value2 = 15200
$ws.Range("D25").Value2 = 0
$ws.Range("E25").Value2 = 0
$ws.Range("D26").Value2 = 111600
$ws.Range("E26").Value2 = 84100
$ws.Range("D27").Value2 = 111600
$ws.Range("E27").Value2 = 84100
$ws.Range("D28").Value2 = 0
$ws.Range("E28").Value2 = 0
$ws.Range("D29").Value2 = 0
$ws.Range("E29").Value2 = 0
$ws.Range("D30").Value2 = 0
$ws.Range("E30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("E31").Value2 = 0
$ws.Range("D32").Value2 = -900
$ws.Range("E32").Value2 = -800
$ws.Range("D33").Value2 = 111600
$ws.Range("E33").Value2 = 84100
$ws.Range("D34").Value2 = 0
$ws.Range("E34").Value2 = 0
$ws.Range("D35").Value2 = 111600
$ws.Range("E35").Value2 = 84100
$ws.Range("D38").Value2 = 43464
$ws.Range("E38").Value2 = 43352
$ws.Range("D41").Value2 = 25400
$ws.Range("E41").Value2 = 84600
$ws.Range("D42").Value2 = 0
$ws.Range("E42").Value2 = 0
$ws.Range("D43").Value2 = 190100
$ws.Range("E43").Value2 = 170200
$ws.Range("D44").Value2 = 46000
$ws.Range("E44").Value2 = 41400
$ws.Range("D45").Value2 = 305400
$ws.Range("E45").Value2 = 309400
$ws.Range("D46").Value2 = 567000
$ws.Range("E46").Value2 = 605600
$ws.Range("D47").Value2 = 8700
$ws.Range("E47").Value2 = "NA"
$ws.Range("D48").Value2 = 234900
$ws.Range("E48").Value2 = 207000
$ws.Range("D49").Value2 = 78700
$ws.Range("E49").Value2 = 74800
$ws.Range("D50").Value2 = 0
$ws.Range("E50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("E51").Value2 = 0
$ws.Range("D52").Value2 = 18000
$ws.Range("E52").Value2 = 24700
$ws.Range("D53").Value2 = 0
$ws.Range("E53").Value2 = 0
$ws.Range("D54").Value2 = 907400
$ws.Range("E54").Value2 = 912100
$ws.Range("D57").Value2 = 92500
$ws.Range("E57").Value2 = 87500
$ws.Range("D58").Value2 = 35900
$ws.Range("E58").Value2 = 35800
$ws.Range("D59").Value2 = 251300
$ws.Range("E59").Value2 = 253100
$ws.Range("D60").Value2 = 379700
$ws.Range("E60").Value2 = 376400
$ws.Range("D61").Value2 = 3495700
$ws.Range("E61").Value2 = 3437600
$ws.Range("D62").Value2 = 71900
$ws.Range("E62").Value2 = 71800
$ws.Range("D63").Value2 = 0
$ws.Range("E63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("E64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("E65").Value2 = 0
$ws.Range("D66").Value2 = 3947300
$ws.Range("E66").Value2 = 3885900
$ws.Range("D68").Value2 = 0
$ws.Range("E68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("E69").Value2 = 0
$ws.Range("D70").Value2 = 0
$ws.Range("E70").Value2 = 0
$ws.Range("D71").Value2 = 0
$ws.Range("E71").Value2 = 0
$ws.Range("D72").Value2 = -3036500
$ws.Range("E72").Value2 = -2972600
$ws.Range("D73").Value2 = 0
$ws.Range("E73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("E74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("E75").Value2 = 0
$ws.Range("D76").Value2 = -3039900
$ws.Range("E76").Value2 = -2973800
$ws.Range("D77").Value2 = 0
$ws.Range("E77").Value2 = 0
$ws.Range("D80").Value2 = 43464
$ws.Range("E80").Value2 = 43352
$ws.Range("D81").Value2 = 111600
$ws.Range("E81").Value2 = 84100
$ws.Range("D83").Value2 = 17900
$ws.Range("E83").Value2 = 12500
$ws.Range("D84").Value2 = 0
$ws.Range("E84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("E85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("E86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("E87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("E88").Value2 = 0
$ws.Range("D89").Value2 = 131700
$ws.Range("E89").Value2 = 107800
$ws.Range("D91").Value2 = -54800
$ws.Range("E91").Value2 = -27800
$ws.Range("D92").Value2 = 0
$ws.Range("E92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("E93").Value2 = 0
$ws.Range("D94").Value2 = -22900
$ws.Range("E94").Value2 = -21600
$ws.Range("D96").Value2 = -45400
$ws.Range("E96").Value2 = -23200
$ws.Range("D97").Value2 = 0
$ws.Range("E97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("E98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("E99").Value2 = 0
$ws.Range("D100").Value2 = -150700
$ws.Range("E100").Value2 = -142000
$ws.Range("D101").Value2 = -300
$ws.Range("E101").Value2 = -100
$ws.Range("D102").Value2 = -42300
$ws.Range("E102").Value2 = -55900

$wb.Save()
